$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New rows of accelerometer data to insert right after the header row (row 1).
# This pushes the existing data rows (old rows 2-19) down by 5.
$newRows = @(
    @(-0.7120251655578613, 1.784507364034653, 2.399995267391204),
    @(-0.8870223164558411, 1.693844005465508, 2.417137637734413),
    @(-1.0569589138031,    1.552494168281555, 2.973462641239166),
    @(-1.252092391252517,  1.304344907402993, 3.757617935538292),
    @(-2.462630152702332,  0.6072362959384919, 3.575700670480729)
)

# Insert 5 blank rows starting at row 2, shifting rows 2..19 down to rows 7..24.
$ws.Range("A2:A6").EntireRow.Insert()

# The insert above copies the header row's (bold/centered) formatting onto the
# new rows; strip it back to the plain/default look used by the other data rows.
$ws.Range("A2:C6").ClearFormats()

# Fill the newly inserted rows with the new data.
for ($i = 0; $i -lt $newRows.Count; $i++) {
    $r = 2 + $i
    $ws.Cells.Item($r, 1).Value = $newRows[$i][0]
    $ws.Cells.Item($r, 2).Value = $newRows[$i][1]
    $ws.Cells.Item($r, 3).Value = $newRows[$i][2]
}

# The old last three data rows (formerly rows 17-19, now rows 22-24) are dropped.
$ws.Range("A22:A24").EntireRow.Delete()
